$d = $word.ActiveDocument

$replacements = @(
    @("35.35 +/- 0.04", "41.48 +/- 0.03"),
    @("5,620 (2.63)", "0 (0.00)"),
    @("10,215 (4.79)", "7,753 (4.97)"),
    @("3,517 (1.65)", "2,535 (1.62)"),
    @("98,718 (46.28)", "75,699 (48.52)"),
    @("113,321 (53.12)", "79,660 (51.05)"),
    @("1,284 (0.60)", "669 (0.43)"),
    @("7,654 (3.59)", "5,835 (3.74)"),
    @("12,165 (5.70)", "9,657 (6.19)"),
    @("29,444 (13.80)", "19,648 (12.59)"),
    @("944 (0.44)", "673 (0.43)"),
    @("7,800 (3.66)", "5,490 (3.52)"),
    @("15,621 (7.32)", "10,122 (6.49)"),
    @("139,695 (65.49)", "104,603 (67.04)"),
    @("213,323 (100.00)", "156,028 (100.00)"),
    @("19,765 (9.27)", "14,952 (9.58)"),
    @("77,090 (36.14)", "60,921 (39.04)"),
    @("116,468 (54.60)", "80,155 (51.37)"),
    @("41,076 (19.26)", "30,286 (19.41)"),
    @("45,665 (21.41)", "33,951 (21.76)"),
    @("70,242 (32.93)", "52,482 (33.64)"),
    @("56,340 (26.41)", "39,309 (25.19)"),
    @("39,034 (18.30)", "29,204 (18.72)"),
    @("47,023 (22.04)", "34,807 (22.31)"),
    @("58,523 (27.43)", "42,947 (27.53)"),
    @("68,743 (32.22)", "49,070 (31.45)"),
    @("811 (0.38)", "600 (0.38)"),
    @("32,153 (15.07)", "21,014 (13.47)"),
    @("108 (0.05)", "64 (0.04)"),
    @("48,239 (22.61)", "37,105 (23.78)"),
    @("131.0 (0.06)", "130 (0.08)"),
    @("1,328.0 (0.62)", "1,252 (0.80)"),
    @("1,264.0 (0.59)", "1,130 (0.72)"),
    @("254.0 (0.12)", "208 (0.13)"),
    @("446.0 (0.21)", "392 (0.25)"),
    @("642.0 (0.30)", "544 (0.35)"),
    @("6,983.0 (3.27)", "6,456 (4.14)"),
    @("8,160.0 (3.83)", "7,114 (4.56)"),
    @("1,167.0 (0.55)", "1,005 (0.64)"),
    @("27,080.0 (12.69)", "24,115 (15.46)"),
    @("24,527.0 (11.50)", "21,827 (13.99)"),
    @("11,362.0 (5.33)", "7,809 (5.00)"),
    @("14,138.0 (6.63)", "12,704 (8.14)"),
    @("801.0 (0.38)", "572 (0.37)"),
    @("6,812.0 (3.19)", "6,171 (3.96)"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done."
